# This script reorders (permutes) rows 5-23 of the active worksheet so that
# each full record (columns A:AY) ends up on the row indicated by the
# mapping below. Row 12 and row 21 remain where they are.
#
# Mapping is "destination row" -> "source row" (i.e. which row's current
# content must end up at the destination row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"
$firstRow = 5
$lastRow = 23

# destination row -> source row
$mapping = @{
    5  = 6
    6  = 7
    7  = 5
    8  = 19
    9  = 10
    10 = 11
    11 = 9
    12 = 12
    13 = 16
    14 = 17
    15 = 14
    16 = 18
    17 = 8
    18 = 13
    19 = 15
    20 = 22
    21 = 21
    22 = 23
    23 = 20
}

# 1) Snapshot every source row's full contents (A:AY) BEFORE writing
#    anything, so overlapping reads/writes of the permutation do not
#    clobber data we still need to read later.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $snapshot[$r] = $rng.Value()
}

# 2) Write each snapshot back out to its destination row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $destRng = $ws.Range("A" + $destRow + ":" + $lastCol + $destRow)
    $destRng.Value = $snapshot[$srcRow]
}
